$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised historical values (rows 259, 262, 263, 265, 266) ---
$ws.Range("B259").Value = 562328000000
$ws.Range("D259").Value = 154472982995.9069

$ws.Range("B262").Value = 557944000000
$ws.Range("D262").Value = 153323440505.6335

$ws.Range("B263").Value = 554911000000
$ws.Range("D263").Value = 152448076923.0769

$ws.Range("B265").Value = 554033000000
$ws.Range("D265").Value = 152202686739.3753

$ws.Range("B266").Value = 551532000000
$ws.Range("D266").Value = 151465685332.1616

# --- Revised latest row (row 313) ---
$ws.Range("B313").Value = 683766000000
$ws.Range("D313").Value = 188502901567.2599

# --- New rows 314 and 315, matching the date-column style from row 313 ---
$ws.Range("A313").Copy()
$ws.Range("A314").PasteSpecial(-4122)
$ws.Range("A315").PasteSpecial(-4122)

$ws.Range("A314").Value = 45139
$ws.Range("B314").Value = 685009000000
$ws.Range("C314").Value = 0.2747003705707999
$ws.Range("D314").Value = 188172226144.3331

$ws.Range("A315").Value = 45170
$ws.Range("B315").Value = 702188000000
$ws.Range("C315").Value = 0.2747252747252747
$ws.Range("D315").Value = 192908791208.7912
